# Chloride field sampling update - 8/11 sampling round
# Appends one new reading row to each field-data worksheet and updates
# the saved cell-selection on every sheet (matching a user who scrolled
# through each tab after pasting in the new data).

$wb = $excel.ActiveWorkbook

# --- WIC -----------------------------------------------------------
$ws = $wb.Worksheets.Item("WIC")
$ws.Range("A4").Value = 44054.365972222222
$ws.Range("B4").Value = 46.46
$ws.Range("C4").Value = 23
[void]$ws.Range("B5").Select()

# --- YS --------------------------------------------------------------
$ws = $wb.Worksheets.Item("YS")
$ws.Range("A17").Value = 44054.487500000003
$ws.Range("B17").Value = 15.52
$ws.Range("C17").Value = 24.5
[void]$ws.Range("A17").Select()

# --- SW ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("SW")
$ws.Range("A17").Value = 44054.505555555559
$ws.Range("B17").Value = 51.96
$ws.Range("C17").Value = 22.6
[void]$ws.Range("A17:C17").Select()

# --- YI ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YI")
$ws.Range("A17").Value = 44054.390972222223
$ws.Range("B17").Value = 13.26
$ws.Range("C17").Value = 24.3
[void]$ws.Range("A17:C17").Select()

# --- YN --------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YN")
$ws.Range("A17").Value = 44054.408333333333
$ws.Range("B17").Value = 13.36
$ws.Range("C17").Value = 24.6
[void]$ws.Range("A17:C17").Select()

# --- 6MC ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("6MC")
$ws.Range("A17").Value = 44054.422222222223
$ws.Range("B17").Value = 13.77
$ws.Range("C17").Value = 21.2
[void]$ws.Range("I20").Select()

# --- DC ------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DC")
$ws.Range("A17").Value = 44054.429861111108
$ws.Range("B17").Value = 13.77
$ws.Range("C17").Value = 17.9
[void]$ws.Range("B17").Select()

# --- PBMS ------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PBMS")
$ws.Range("A17").Value = 44054.461111111108
$ws.Range("B17").Value = 46.87
$ws.Range("C17").Value = 21.9
[void]$ws.Range("A17:C17").Select()

# --- PBSF --------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PBSF")
$ws.Range("A17").Value = 44054.45
$ws.Range("B17").Value = 223.08
$ws.Range("C17").Value = 22.9
[void]$ws.Range("H13").Select()

# --- MO (no new data this round, just a moved selection; stays the active tab) -----
$ws = $wb.Worksheets.Item("MO")
[void]$ws.Range("F19").Select()
